# BUS_003 resume & Big Interview
#
# The elevator-pitch paragraph's closing sentence is expanded: the old
# "...research assistant positions during which exemplify analysis and
# technical skills." tail is replaced with a longer, corrected passage
# describing the research-assistant position, the Leader Bank internship,
# and a closing aspirational sentence.

$d = $word.ActiveDocument

$oldTail = "research assistant positions during which exemplify analysis and technical skills."

$newTail = "a research assistant position which have allowed me to exemplify my leadership and organizational skills. Most recently I have interned at Leader Bank. This prepared me for efficient operation in the corporate financial landscape and inspired me to pursue a career in the Finical Technologies (FinTech) industry. " + `
    "I hope to use the skills I" + [char]0x2019 + "ve honed with whichever employs me in the future, and I intend to move up corporate ladders and make a difference at any organization."

$found = $d.Content.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2)

if (-not $found) {
    throw "Could not locate the elevator-pitch closing sentence to replace."
}
